# edit.ps1 - Apply the "Added NicsView..." commit's changes:
#   1. Refresh the cached datetimeFigureOut field text (7/26/2024 -> 11/2/2024)
#      on the slide master and every slide layout's Date placeholder.
#   2. Reposition/resize nine "Rounded Rectangle" icon shapes on slide 8
#      (the Hexulator icon) to make room for new fields.

$p = $ppt.ActivePresentation
$EMU_PER_PT = 12700
# Shape.Top/.Left/.Width/.Height are COM `Single` (float32) properties, so a
# plain EMU/12700 conversion can land a hair under the target EMU once it
# round-trips through float32 and back to integer EMU on save. A tiny nudge
# (well under half an EMU) compensates without risking overshoot.
$EMU_EPS_PT = 0.00003

# ---------------------------------------------------------------------------
# 1) Date placeholder field text: 7/26/2024 -> 11/2/2024
#    (slide master + all 13 slide layouts)
# ---------------------------------------------------------------------------
function Update-DateField {
    param($shapes)
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "7/26/2024") {
                $sh.TextFrame.TextRange.Text = "11/2/2024"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    Update-DateField $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Reposition / resize the Hexulator icon shapes on slide 8
# ---------------------------------------------------------------------------
$slide8 = $p.Slides.Item(8)

# shapePos -> @{ Top = newTopEmu; Height = newHeightEmu (optional) }
$moves = @(
    @{ Pos = 12; Top = 4023684 },                       # Rectangle: Rounded Corners 34
    @{ Pos = 13; Top = 4023684; Height = 566738 },       # Rectangle: Rounded Corners 35
    @{ Pos = 14; Top = 4365384 },                        # Rectangle: Rounded Corners 36
    @{ Pos = 15; Top = 4361946 },                        # Rectangle: Rounded Corners 37
    @{ Pos = 16; Top = 4685903 },                        # Rectangle: Rounded Corners 38
    @{ Pos = 17; Top = 4365389 },                        # Rectangle: Rounded Corners 39
    @{ Pos = 18; Top = 4380938; Height = 566738 },       # Rectangle: Rounded Corners 40
    @{ Pos = 19; Top = 4020870; Height = 650358 },       # Rectangle: Rounded Corners 41
    @{ Pos = 20; Top = 4020557 }                         # Rectangle: Rounded Corners 42
)

foreach ($m in $moves) {
    $sh = $slide8.Shapes.Item($m.Pos)
    $sh.Top = ($m.Top / $EMU_PER_PT) + $EMU_EPS_PT
    if ($m.ContainsKey("Height")) {
        $sh.Height = ($m.Height / $EMU_PER_PT) + $EMU_EPS_PT
    }
}
